$d = $word.ActiveDocument

$pairs = @(
  @("337÷2=", "281÷8="),
  @("268÷7=", "689÷8="),
  @("810÷4=", "676÷2="),
  @("746÷2=", "746÷6="),
  @("129÷9=", "998÷2="),
  @("346÷2=", "893÷8="),
  @("331÷7=", "150÷9="),
  @("232÷9=", "863÷6="),
  @("956÷8=", "992÷7="),
  @("404÷2=", "973÷4="),
  @("922÷6=", "731÷7="),
  @("768÷4=", "197÷2="),
  @("302÷8=", "311÷3="),
  @("782÷6=", "428÷9="),
  @("939÷8=", "884÷2="),
  @("261÷8=", "523÷4="),
  @("209÷5=", "843÷6="),
  @("318÷5=", "516÷8="),
  @("844÷3=", "319÷4="),
  @("881÷9=", "247÷8="),
  @("444÷9=", "134÷2="),
  @("699÷6=", "685÷5="),
  @("744÷7=", "220÷5="),
  @("103÷2=", "600÷6="),
  @("176÷4=", "615÷3=")
)

foreach ($pair in $pairs) {
  $old = $pair[0]
  $new = $pair[1]
  $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                           $true, 1, $false, $new, 2)
}
